$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.17%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.77%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.263"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.07%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08109"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.30%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.523"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.09%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.638"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.33%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.911"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.68%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.942"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.54%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9346"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.61%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1333"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "18.79%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1951"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.80%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09200"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.66%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03435"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "4.06%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09542"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.64%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001395"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.56%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006021"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.07%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.359"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.59%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.36%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.226"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "22.00%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1313"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.63%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2311"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.73%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.51%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001222"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.17%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004359"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.00%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.16%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003991"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.04%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02493"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "11.30%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05246"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.96%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007685"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.66%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1431"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.65%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008603"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.55%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002160"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.39%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008170"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.43%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006660"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.30%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002852"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-13.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.05%"
